$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert Q3:V3 into a shared formula group (matches the target workbook's
# restructuring of the Q3:V3 formulas into a single shared formula).
$ws.Range("Q3:V3").Formula = '=Q2*$A$29'

# Add the new "Items:" / "Type:" / "Equipment ID" table in rows 34-36.
$ws.Range("A34").Value = "Items:"
$ws.Range("A35").Value = "Type:"
$ws.Range("A36").Value = "Equipment ID"

$ws.Range("B35").Value = "Sword"
$ws.Range("E35").Value = "Hat"
$ws.Range("H35").Value = "Boots"
$ws.Range("G35").Value = "Leggings"
$ws.Range("C35").Value = "Gloves"
$ws.Range("D35").Value = "Secondary"
$ws.Range("F35").Value = "Body"
$ws.Range("I35").Value = "Ring"

$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 3
$ws.Range("F36").Value = 4
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 6
$ws.Range("I36").Value = 7

# Update the view state to match: scroll down so row 22 is at top, and
# select cell A37 (just below the new table).
$ws.Activate()
$null = $ws.Range("A37").Select()
